$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 863-864, shifting existing rows 863.. down by 2.
$ws.Rows("863:864").Insert()

# New row 863
$ws.Cells.Item(863, 1).Value2 = 5
$ws.Cells.Item(863, 2).Value2 = "Macroferia Regional de Talca"
$ws.Cells.Item(863, 3).Value2 = "Maule"
$ws.Cells.Item(863, 4).Value2 = 45106
$ws.Cells.Item(863, 5).Value2 = 7
$ws.Cells.Item(863, 6).Value2 = 100112004
$ws.Cells.Item(863, 7).Value2 = "Cebolla"
$ws.Cells.Item(863, 8).Value2 = "Sin especificar"
$ws.Cells.Item(863, 9).Value2 = "1a (guarda)"
$ws.Cells.Item(863, 10).Value2 = 2000
$ws.Cells.Item(863, 11).Value2 = 10000
$ws.Cells.Item(863, 12).Value2 = 10000
$ws.Cells.Item(863, 13).Value2 = 10000
$ws.Cells.Item(863, 14).Value2 = "$/malla 25 kilos"
$ws.Cells.Item(863, 15).Value2 = "Región del Maule"
$ws.Cells.Item(863, 16).Value2 = 400
$ws.Cells.Item(863, 17).Value2 = 25
$ws.Cells.Item(863, 18).Value2 = "Hortaliza"

# New row 864
$ws.Cells.Item(864, 1).Value2 = 5
$ws.Cells.Item(864, 2).Value2 = "Macroferia Regional de Talca"
$ws.Cells.Item(864, 3).Value2 = "Maule"
$ws.Cells.Item(864, 4).Value2 = 45106
$ws.Cells.Item(864, 5).Value2 = 7
$ws.Cells.Item(864, 6).Value2 = 100112004
$ws.Cells.Item(864, 7).Value2 = "Cebolla"
$ws.Cells.Item(864, 8).Value2 = "Sin especificar"
$ws.Cells.Item(864, 9).Value2 = "2a (guarda)"
$ws.Cells.Item(864, 10).Value2 = 600
$ws.Cells.Item(864, 11).Value2 = 9000
$ws.Cells.Item(864, 12).Value2 = 9000
$ws.Cells.Item(864, 13).Value2 = 9000
$ws.Cells.Item(864, 14).Value2 = "$/malla 25 kilos"
$ws.Cells.Item(864, 15).Value2 = "Región del Maule"
$ws.Cells.Item(864, 16).Value2 = 360
$ws.Cells.Item(864, 17).Value2 = 25
$ws.Cells.Item(864, 18).Value2 = "Hortaliza"

# Ensure date formatting (numFmt used for column D, style index 2 in the original file)
$ws.Range("D863:D864").NumberFormat = "YYYY-MM-DD HH:MM:SS"
